$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
  @(1108, 475, "V", "Chicago", 0, 3, 3, 21, 27, 40.5, 40, 260, 18.5),
  @(1108, 476, "H", "Pittsburgh", 7, 7, 6, 9, 29, -5, 7.5, -310, -3),
  @(1111, 113, "V", "Baltimore", 3, 0, 0, 7, 10, 7, 10, -410, 6.5),
  @(1111, 114, "H", "Miami", 0, 6, 0, 16, 22, 48, 47, 340, 22.5),
  @(1114, 241, "V", "Jacksonville", 6, 3, 0, 8, 17, 48, 48, 350, 23),
  @(1114, 242, "H", "Indianapolis", 17, 3, 0, 3, 23, -10.5, 11.5, -420, 3.5),
  @(1114, 243, "V", "Cleveland", 7, 0, 0, 0, 7, 44.5, 44.5, 125, 0.5),
  @(1114, 244, "H", "NewEngland", 7, 17, 7, 14, 45, -3, 2.5, -145, 21.5),
  @(1114, 245, "V", "Atlanta", 3, 0, 0, 0, 3, 51.5, 55.5, 300, -0.5),
  @(1114, 246, "H", "Dallas", 7, 29, 7, 0, 43, -9.5, 9, -360, 24),
  @(1114, 247, "V", "Buffalo", 10, 7, 21, 7, 45, -13.5, 14, -800, 6.5),
  @(1114, 248, "H", "NYJets", 0, 3, 0, 14, 17, 48.5, 48.5, 550, 21.5),
  @(1114, 249, "V", "NewOrleans", 0, 6, 6, 9, 21, 44.5, 42.5, 130, -0.5),
  @(1114, 250, "H", "Tennessee", 3, 10, 7, 3, 23, -2.5, 3, -150, 21),
  @(1114, 251, "V", "TampaBay", 0, 6, 7, 6, 19, -8, 10.5, -430, -7.5),
  @(1114, 252, "H", "Washington", 6, 10, 7, 6, 29, 51, 51, 360, 27),
  @(1114, 253, "V", "Detroit", 0, 10, 6, 0, 16, 44, 40.5, 210, 20),
  @(1114, 254, "H", "Pittsburgh", 7, 3, 3, 3, 16, -9, -6, -250, 3),
  @(1114, 255, "V", "Minnesota", 0, 13, 7, 7, 27, 50.5, 53.5, 155, 26.5),
  @(1114, 256, "H", "LAChargers", 0, 10, 7, 3, 20, 3, 3.5, -175, 4),
  @(1114, 257, "V", "Carolina", 17, 6, 8, 3, 34, 45.5, 41.5, 270, 20),
  @(1114, 258, "H", "Arizona", 0, 0, 3, 7, 10, 11.5, 7.5, -330, 3.5),
  @(1114, 259, "V", "Seattle", 0, 0, 0, 0, 0, 49.5, 49, 145, 24),
  @(1114, 260, "H", "GreenBay", 0, 3, 0, 14, 17, -6, 3, -165, 0.5),
  @(1114, 261, "V", "Philadelphia", 10, 10, 7, 3, 30, 44.5, 44.5, -105, 21),
  @(1114, 262, "H", "Denver", 0, 10, 3, 0, 13, -1.5, "pk", -115, -3),
  @(1114, 263, "V", "KansasCity", 7, 10, 10, 14, 41, -3, 2.5, -145, 27),
  @(1114, 264, "H", "LasVegas", 0, 7, 7, 0, 14, 53.5, 53.5, 125, 3),
  @(1115, 265, "V", "LARams", 0, 7, 0, 3, 10, 3, 3.5, -190, -6),
  @(1115, 266, "H", "SanFrancisco", 14, 7, 3, 7, 31, 47, 50, 170, 27)
)

$startRow = 272
$nRows = $rowsData.Count
$nCols = 13
$arr = New-Object 'object[,]' $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
  for ($c = 0; $c -lt $nCols; $c++) {
    $arr[$r, $c] = $rowsData[$r][$c]
  }
}

$endRow = $startRow + $nRows - 1
$rng = $ws.Range("A$startRow" + ":M$endRow")
$rng.Value = $arr
